# Fill in the newly-computed "CoT" (chain-of-thought) results into the
# "llm (values)" worksheet. The "Zero-Shot (CoT)" (rows 21-27) and
# "Few-Shot (CoT)" (rows 35-41) blocks already existed as formatted, but
# empty, placeholder rows; this adds the measured accuracy values for the
# Hint/Value/Comb columns (I, J, L) for each heuristic (angle, dist, area,
# comb a-d, comb a-m, comb d-m, all).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("llm (values)")

# --- Zero-Shot (CoT) ------------------------------------------------------
$ws.Range("I21").Value = 0.68200000000000005
$ws.Range("J21").Value = 0.94199999999999995
$ws.Range("L21").Value = 0.94699999999999995

$ws.Range("I22").Value = 0.59899999999999998
$ws.Range("J22").Value = 0.61599999999999999
$ws.Range("L22").Value = 0.74199999999999999

$ws.Range("I23").Value = 0.59899999999999998
$ws.Range("J23").Value = 0.60199999999999998
$ws.Range("L23").Value = 0.64300000000000002

$ws.Range("I24").Value = 0.63400000000000001
$ws.Range("J24").Value = 0.94599999999999995
$ws.Range("L24").Value = 0.95

$ws.Range("I25").Value = 0.59899999999999998
$ws.Range("J25").Value = 0.92
$ws.Range("L25").Value = 0.96299999999999997

$ws.Range("I26").Value = 0.60699999999999998
$ws.Range("J26").Value = 0.73199999999999998
$ws.Range("L26").Value = 0.80700000000000005

$ws.Range("I27").Value = 0.60299999999999998
$ws.Range("J27").Value = 0.91400000000000003
$ws.Range("L27").Value = 0.95399999999999996

# --- Few-Shot (CoT) --------------------------------------------------------
$ws.Range("I35").Value = 0.93899999999999995
$ws.Range("J35").Value = 0.92200000000000004
$ws.Range("L35").Value = 0.93700000000000006

$ws.Range("I36").Value = 0.54
$ws.Range("J36").Value = 0.80700000000000005
$ws.Range("L36").Value = 0.84

$ws.Range("I37").Value = 0.66900000000000004
$ws.Range("J37").Value = 0.749
$ws.Range("L37").Value = 0.71699999999999997

$ws.Range("I38").Value = 0.57699999999999996
$ws.Range("J38").Value = 0.94299999999999995
$ws.Range("L38").Value = 0.94399999999999995

# NOTE: this one is recorded as the literal text "0..610" (a typo for
# "0.610" made when the value was transcribed), not as a number.
$ws.Range("I39").Value = "0..610"
$ws.Range("J39").Value = 0.95499999999999996
$ws.Range("L39").Value = 0.97099999999999997

$ws.Range("I40").Value = 0.34100000000000003
$ws.Range("J40").Value = 0.89800000000000002
$ws.Range("L40").Value = 0.92100000000000004

$ws.Range("I41").Value = 0.60899999999999999
$ws.Range("J41").Value = 0.96499999999999997
$ws.Range("L41").Value = 0.97699999999999998

# --- cosmetic: leave the selection/active cell where the author last
# clicked while reviewing the other sheets.
$ws2 = $wb.Worksheets.Item("llm (hint)")
$ws2.Range("P26").Select()

$ws.Activate()
$ws.Range("N30").Select()
